# Re-colour the presentation's live theme (ppt/theme/theme2.xml - the part actually
# wired to the slide master / presentation via the `theme` relationship) from the
# "Integral" / "Red Violet" palette back to the stock PowerPoint "Office" palette.
#
# PowerPoint COM exposes the 12 theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) as a 1-based collection on ThemeColorScheme, reachable from any slide.
# RGB values use the standard VBA BGR-packed integer encoding: R + G*256 + B*65536.

function BGR($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Index : Name     : Office Theme hex
#   1   : dk1      : 000000
#   2   : lt1      : FFFFFF
#   3   : dk2      : 44546A
#   4   : lt2      : E7E6E6
#   5   : accent1  : 5B9BD5
#   6   : accent2  : ED7D31
#   7   : accent3  : A5A5A5
#   8   : accent4  : FFC000
#   9   : accent5  : 4472C4
#  10   : accent6  : 70AD47
#  11   : hlink    : 0563C1
#  12   : folHlink : 954F72

$tcs.Colors(1).RGB  = BGR 0x00 0x00 0x00
$tcs.Colors(2).RGB  = BGR 0xFF 0xFF 0xFF
$tcs.Colors(3).RGB  = BGR 0x44 0x54 0x6A
$tcs.Colors(4).RGB  = BGR 0xE7 0xE6 0xE6
$tcs.Colors(5).RGB  = BGR 0x5B 0x9B 0xD5
$tcs.Colors(6).RGB  = BGR 0xED 0x7D 0x31
$tcs.Colors(7).RGB  = BGR 0xA5 0xA5 0xA5
$tcs.Colors(8).RGB  = BGR 0xFF 0xC0 0x00
$tcs.Colors(9).RGB  = BGR 0x44 0x72 0xC4
$tcs.Colors(10).RGB = BGR 0x70 0xAD 0x47
$tcs.Colors(11).RGB = BGR 0x05 0x63 0xC1
$tcs.Colors(12).RGB = BGR 0x95 0x4F 0x72
